# Apply the "Add files via upload" revision to the USER STORIES - CAMPUS
# RECRUITMENT SYSTEM workbook.
#
# The edit restructures the single worksheet: it adds a second task column
# (TASK1 / TASK2), merges column A per-section (STUDENT / TPO / COMPANY),
# adds a "Generate Result" task under the TPO > View Student branch, adds a
# "Send notification to TPO" task under COMPANY, fixes a couple of typos,
# and moves a few labels around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Start from a clean merge state so the new merged ranges below don't
#    coalesce with the old ones.
# ---------------------------------------------------------------------
$ws.Cells.UnMerge()

# ---------------------------------------------------------------------
# 2. Clear the couple of cells whose content moved / disappeared and is
#    not going to be overwritten by the same-address assignments below.
# ---------------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = ""   # old "COMPANY" label (now at A23)
$ws.Cells.Item(25, 2).Value = ""   # old "Rrgistration" typo (removed)

# ---------------------------------------------------------------------
# 3. Write the final cell values (row/column absolute addressing so the
#    row-shift caused by the two newly inserted rows is handled simply
#    by writing the destination state directly).
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "USER STORIES - CAMPUS RECRUITMENT SYSTEM"

$ws.Cells.Item(2, 1).Value = "FEATURES"
$ws.Cells.Item(2, 2).Value = "SUB-FEATURES"
$ws.Cells.Item(2, 3).Value = "TASK1"
$ws.Cells.Item(2, 4).Value = "TASK2"

# STUDENT section (rows 3-11)
$ws.Cells.Item(3, 1).Value = "STUDENT"
$ws.Cells.Item(3, 2).Value = "Login"
$ws.Cells.Item(4, 3).Value = "Forgot Password"
$ws.Cells.Item(5, 3).Value = "Reset Password"
$ws.Cells.Item(6, 2).Value = "Registration"
$ws.Cells.Item(6, 3).Value = "Create Profile"
$ws.Cells.Item(7, 3).Value = "Update Profile"
$ws.Cells.Item(8, 2).Value = "View Company"
$ws.Cells.Item(8, 3).Value = "Search by Salary"
$ws.Cells.Item(9, 3).Value = "Search by Job"
$ws.Cells.Item(10, 3).Value = "Apply for Drive"
$ws.Cells.Item(11, 2).Value = "                       View Materials"
$ws.Cells.Item(11, 3).Value = "Download and Read"

# TPO section (rows 12-22)
$ws.Cells.Item(12, 1).Value = "TPO"
$ws.Cells.Item(12, 2).Value = "Login"
$ws.Cells.Item(13, 3).Value = "Forgot Password"
$ws.Cells.Item(14, 3).Value = "Reset Password"
$ws.Cells.Item(15, 2).Value = "                          Registration"
$ws.Cells.Item(16, 2).Value = "View Company"
$ws.Cells.Item(16, 3).Value = "Accept"
$ws.Cells.Item(17, 3).Value = "Reject"
$ws.Cells.Item(18, 3).Value = "View"
$ws.Cells.Item(19, 2).Value = "View Student"
$ws.Cells.Item(19, 3).Value = "Accept"
$ws.Cells.Item(19, 4).Value = "Send notification"
$ws.Cells.Item(20, 3).Value = "Reject"
$ws.Cells.Item(21, 2).Value = "View Status"
$ws.Cells.Item(21, 3).Value = "Get Result"
$ws.Cells.Item(22, 2).Value = "                          Materials"
$ws.Range("B22").Style = "Normal"   # drop the stale centered style this cell had as the old B22 "login"
$ws.Cells.Item(22, 3).Value = "Upload materials"

# COMPANY section (rows 23-29)
$ws.Cells.Item(23, 1).Value = "COMPANY"
$ws.Cells.Item(23, 2).Value = "login"
$ws.Cells.Item(23, 3).Value = "Forgot Password"
$ws.Cells.Item(24, 3).Value = "Reset Password"
$ws.Cells.Item(25, 3).Value = "Create Profile"
$ws.Cells.Item(26, 2).Value = "Registration"
$ws.Cells.Item(26, 3).Value = "Post Job Vacancy"
$ws.Cells.Item(27, 3).Value = "Qualification & Criterias"
$ws.Cells.Item(28, 2).Value = "                     Generate Result"
$ws.Cells.Item(28, 3).Value = "Send to TPO"
$ws.Cells.Item(29, 2).Value = "              Send notification to TPO"

# ---------------------------------------------------------------------
# 4. Match the bold 14pt header formatting (style used by A2:C2) on the
#    newly added D2 header cell.
# ---------------------------------------------------------------------
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Font.Size = 14

# ---------------------------------------------------------------------
# 5. Center (horizontal + vertical) the section/sub-feature columns.
#    Column A spans the whole STUDENT/TPO/COMPANY block (rows 3-29);
#    column B gets the same treatment except for the narrower
#    "indented" sub-rows (11, 15, 21, 22, 28, 29) which keep their own
#    look.
# ---------------------------------------------------------------------
$ws.Range("A3:A29").VerticalAlignment = -4108
$ws.Range("A3:A29").HorizontalAlignment = -4108

$ws.Range("B3:B10").VerticalAlignment = -4108
$ws.Range("B3:B10").HorizontalAlignment = -4108

$ws.Range("B12:B14").VerticalAlignment = -4108
$ws.Range("B12:B14").HorizontalAlignment = -4108

$ws.Range("B16:B20").VerticalAlignment = -4108
$ws.Range("B16:B20").HorizontalAlignment = -4108

$ws.Range("B21").VerticalAlignment = -4108
$ws.Range("B21").HorizontalAlignment = -4108

$ws.Range("B23:B27").VerticalAlignment = -4108
$ws.Range("B23:B27").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Re-create the merged ranges for the new layout.
# ---------------------------------------------------------------------
$ws.Range("A3:A11").Merge()
$ws.Range("A12:A22").Merge()
$ws.Range("A23:A29").Merge()

$ws.Range("B3:B5").Merge()
$ws.Range("B6:B7").Merge()
$ws.Range("B8:B10").Merge()
$ws.Range("B12:B14").Merge()
$ws.Range("B16:B18").Merge()
$ws.Range("B19:B20").Merge()
$ws.Range("B23:B25").Merge()
$ws.Range("B26:B27").Merge()

# ---------------------------------------------------------------------
# 7. Restore the view/selection state recorded in the saved workbook.
# ---------------------------------------------------------------------
$ws.Range("D13").Select()
